$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1 - mirror style of existing header row (e.g. E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style from an existing header cell so the new ones match
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean values for rows 2-4
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = $true

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false
